# Updated legacy GSC export data.
#
# The Google Search Console export window rolled forward by one day:
# the oldest date row ("2025-08-24") drops off, and every subsequent
# day's metrics move up one row to take the place of the previous day.
# Deleting the second row (the first data row, right below the header)
# on the "Chart" sheet reproduces exactly that shift: Excel removes the
# row, pulls every following row up by one, drops the now-unused
# "2025-08-24" shared string, and keeps everything else (the other
# sheets' shared-string references, styles, etc.) internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
